{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Locate the three \"CORE COMPETENCIES\" detail paragraphs ---------------\nlet researchPara = null;\nlet programmingPara = null;\nlet infraPara = null;\n\nfor (const p of items) {\n  const t = p.text;\n  if (t.indexOf(\"Research and Analytics: Survey Methodology:\") === 0) {\n    researchPara = p;\n  } else if (t.indexOf(\"Programming and Development: Python:\") === 0) {\n    programmingPara = p;\n  } else if (t.indexOf(\"Data Infrastructure: Cloud Platforms:\") === 0) {\n    infraPara = p;\n  }\n}\n\nif (researchPara && programmingPara && infraPara) {\n  // Collapse the three detailed bullet paragraphs into a single summary\n  // paragraph, then remove the other two paragraphs entirely.\n  researchPara.insertText(\n    \"Research and Analytics \u2022 Programming and Development \u2022 Data Infrastructure\",\n    \"Replace\"\n  );\n  programmingPara.delete();\n  infraPara.delete();\n  await context.sync();\n}\n\n// --- Insert the new \"TECHNICAL SKILLS\" section -----------------------------\n// Re-load paragraphs fresh (indices/anchors shifted after the deletes above).\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet anchorPara = null;\nfor (const p of paragraphs2.items) {\n  if (\n    p.text.indexOf(\n      \"Built comprehensive survey operations platform from RFP through deployment\"\n    ) !== -1\n  ) {\n    anchorPara = p;\n  }\n}\n\nif (anchorPara) {\n  // Insert each new paragraph empty first, set its style, THEN fill in the\n  // text. Setting .style on a paragraph that already holds text stamps an\n  // explicit (and here redundant, since it equals the inherited default)\n  // <w:pStyle> on it; assigning the style before any text exists keeps the\n  // normal-style paragraphs free of that redundant tag, matching how Word\n  // itself emits them.\n  const heading = anchorPara.insertParagraph(\"\", \"After\");\n  heading.style = \"Heading 2\";\n  heading.insertText(\"TECHNICAL SKILLS\", \"Replace\");\n\n  const p1 = heading.insertParagraph(\"\", \"After\");\n  p1.style = \"Normal\";\n  p1.insertText(\n    \"RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization\",\n    \"Replace\"\n  );\n\n  const p2 = p1.insertParagraph(\"\", \"After\");\n  p2.style = \"Normal\";\n  p2.insertText(\n    \"PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages\",\n    \"Replace\"\n  );\n\n  const p3 = p2.insertParagraph(\"\", \"After\");\n  p3.style = \"Normal\";\n  p3.insertText(\n    \"DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial\",\n    \"Replace\"\n  );\n\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# Part 1: Collapse the three \"CORE COMPETENCIES\" detail paragraphs into a\n# single summary paragraph.\n# ---------------------------------------------------------------------------\n$researchPara = $null\n$programmingPara = $null\n$infraPara = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"Research and Analytics: Survey Methodology:*\") {\n        $researchPara = $p\n    } elseif ($t -like \"Programming and Development: Python:*\") {\n        $programmingPara = $p\n    } elseif ($t -like \"Data Infrastructure: Cloud Platforms:*\") {\n        $infraPara = $p\n    }\n}\n\nif ($researchPara -ne $null -and $programmingPara -ne $null -and $infraPara -ne $null) {\n    $researchPara.Range.Text = \"Research and Analytics \" + [char]0x2022 + \" Programming and Development \" + [char]0x2022 + \" Data Infrastructure\"\n    # Delete the trailing two paragraphs LAST-TO-FIRST. A paragraph object\n    # captured before a structural edit (insert/delete) does not track the\n    # shifted document - deleting $programmingPara first would leave\n    # $infraPara pointing at whatever paragraph landed in its old slot.\n    # Deleting bottom-up means every reference still pending a delete sits\n    # above the edit point, so it stays valid.\n    $infraPara.Range.Delete()\n    $programmingPara.Range.Delete()\n}\n\n# ---------------------------------------------------------------------------\n# Part 2: Insert a new \"TECHNICAL SKILLS\" section after the \"Built\n# comprehensive survey operations platform...\" bullet, before the closing\n# \"For a more detailed...\" paragraph.\n# ---------------------------------------------------------------------------\n$anchorPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Built comprehensive survey operations platform from RFP through deployment*\") {\n        $anchorPara = $p\n    }\n}\n\nif ($anchorPara -ne $null) {\n    # Insert a new empty paragraph after the anchor, set its style BEFORE\n    # giving it text (mirrors Word's own behaviour: a style applied while a\n    # paragraph is still empty stays \"pending\" and whichever style equals the\n    # paragraph's already-inherited default - Normal - never gets stamped as\n    # an explicit <w:pStyle>).\n    $anchorPara.Range.InsertParagraphAfter()\n    $heading = $anchorPara.Next()\n    $heading.Style = \"Heading 2\"\n    $heading.Range.Text = \"TECHNICAL SKILLS\"\n\n    $heading.Range.InsertParagraphAfter()\n    $p1 = $heading.Next()\n    $p1.Style = \"Normal\"\n    $p1.Range.Text = \"RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization\"\n\n    $p1.Range.InsertParagraphAfter()\n    $p2 = $p1.Next()\n    $p2.Style = \"Normal\"\n    $p2.Range.Text = \"PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages\"\n\n    $p2.Range.InsertParagraphAfter()\n    $p3 = $p2.Next()\n    $p3.Style = \"Normal\"\n    $p3.Range.Text = \"DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial\"\n}\n"}
